$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "93.420.87"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "3.115.44"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("D5").Value = "'237.35"
$ws.Range("E5").Value = "  -3.53%  "
$ws.Range("D6").Value = "'613.83"
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "'1.12"
$ws.Range("E7").Value = "  +1.22%  "
$ws.Range("E8").Value = "  +1.82%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'0.841"
$ws.Range("E10").Value = "  +13.73%  "
$ws.Range("D11").Value = "3.114.29"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "'0.198"
$ws.Range("E12").Value = "  -2.62%  "
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").Value = "'35.22"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "93.134.83"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").Value = "'5.43"
$ws.Range("D17").Value = "3.693.04"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "3.106.17"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").Value = "'3.80"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").Value = "'14.86"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("E21").Value = "  +4.31%  "
$ws.Range("D22").Value = "'442.81"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").Value = "'0.0000200"
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("D24").Value = "'9.08"
$ws.Range("E24").Value = "  -4.01%  "
$ws.Range("D25").Value = "'8.23"
$ws.Range("E25").Value = "  +4.86%  "
$ws.Range("D26").Value = "'5.80"
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").Value = "'12.92"
$ws.Range("E27").Value = "  +10.34%  "
$ws.Range("D28").Value = "'85.79"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'0.182"
$ws.Range("E30").Value = "  +9.05%  "
$ws.Range("D31").Value = "'0.241"
$ws.Range("E31").Value = "  +3.56%  "
$ws.Range("E32").Value = "  -13.88%  "
$ws.Range("D33").Value = "'9.26"
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("D34").Value = "'1.01"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").Value = "'8.00"
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("D36").Value = "'0.160"
$ws.Range("E36").Value = "  -9.38%  "
$ws.Range("D37").Value = "'25.97"
$ws.Range("E37").Value = "  -1.09%  "
$ws.Range("D38").Value = "'3.93"
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("D40").Value = "'0.446"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").Value = "'24.00"
$ws.Range("E42").Value = "  +8.07%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'476.87"
$ws.Range("E43").Value = "  -2.97%  "
$ws.Range("D44").Value = "'3.31"
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'159.00"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D47").Value = "'0.702"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("E48").Value = "  -2.24%  "
$ws.Range("D49").Value = "'1.33"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  -0.33%  "
